$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(37).Copy()
$ws.Rows(61).Insert()
